# Insert a new "Industry" column between "Stock Name" (B) and "Mutual Fund" (C),
# shifting the old C:I columns to D:J, then populate the new column with the
# industry classification produced by the motilal_portfolio_change_engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:I one place to the right to make room for the new Industry column.
$ws.Columns.Item(3).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 3).Value = "Industry"

# Industry value per data row (rows 2-26), in the same order as the sheet.
$industries = @(
    "Power",
    "Metals & Minerals Trading",
    "Industrial Products",
    "Power",
    "Minerals & Mining",
    "Cement & Cement Products",
    "Cement & Cement Products",
    "Industrial Products",
    "Ferrous Metals",
    "Fertilizers & Agrochemicals",
    "Banks",
    "Chemicals & Petrochemicals",
    "Power",
    "Power",
    "Power",
    "Fertilizers & Agrochemicals",
    "IT - Software",
    "Non - Ferrous Metals",
    "Fertilizers & Agrochemicals",
    "Electrical Equipment",
    "Power",
    "Oil",
    "Electrical Equipment",
    "Gas",
    "Finance"
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
